$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) - numeric-looking text values; apostrophe-prefixed
# to force Excel to store them as text (matching the source inlineStr type)
# rather than silently parsing them into floating point numbers.
$ws.Range("D2").Value = "'40.791.72"
$ws.Range("D3").Value = "'2.216.41"
$ws.Range("D5").Value = "'229.21"
$ws.Range("D6").Value = "'0.635"
$ws.Range("D7").Value = "'64.47"
$ws.Range("D9").Value = "'0.407"
$ws.Range("D10").Value = "'0.0872"
$ws.Range("D12").Value = "'2.544.06"
$ws.Range("D13").Value = "'15.90"
$ws.Range("D14").Value = "'22.29"
$ws.Range("D15").Value = "'0.823"
$ws.Range("D16").Value = "'5.62"
$ws.Range("D17").Value = "'2.211.07"
$ws.Range("D18").Value = "'40.662.70"
$ws.Range("D19").Value = "'73.88"
$ws.Range("D22").Value = "'253.03"
$ws.Range("D24").Value = "'2.38"
$ws.Range("D27").Value = "'173.29"
$ws.Range("D29").Value = "'20.42"
$ws.Range("D31").Value = "'2.82"
$ws.Range("D32").Value = "'0.124"
$ws.Range("D34").Value = "'7.16"
$ws.Range("D35").Value = "'4.78"
$ws.Range("D36").Value = "'0.0632"
$ws.Range("D37").Value = "'3.84"
$ws.Range("D38").Value = "'2.47"
$ws.Range("D39").Value = "'0.998"
$ws.Range("D40").Value = "'4.88"
$ws.Range("D41").Value = "'8.69"
$ws.Range("D43").Value = "'101.37"
$ws.Range("D44").Value = "'1.23"
$ws.Range("D45").Value = "'1.520.94"
$ws.Range("D46").Value = "'17.36"
$ws.Range("D47").Value = "'0.0940"
$ws.Range("D48").Value = "'1.12"
$ws.Range("D51").Value = "'51.05"

# Row 51: coin name + link text changed (Celestia -> MultiversX)
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"

# Column E (Volume 1h) - percentage text, already non-numeric due to
# leading/trailing spaces and the trailing '%' sign, so plain assignment
# keeps it text without needing the apostrophe prefix.
$ws.Range("E2").Value = "  +3.69%  "
$ws.Range("E3").Value = "  +2.58%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("E18").Value = "  +3.52%  "
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("E20").Value = "  +6.60%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("E22").Value = "  +9.44%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("E25").Value = "  -8.17%  "
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("E28").Value = "  +2.33%  "
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("E30").Value = "  +3.13%  "
$ws.Range("E31").Value = "  +3.22%  "
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  +2.35%  "
$ws.Range("E37").Value = "  +7.44%  "
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").Value = "  +12.21%  "
$ws.Range("E41").Value = "  +11.13%  "
$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("E44").Value = "  +5.05%  "
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("E49").Value = "  +40.20%  "
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("E51").Value = "  +10.23%  "
